$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 58874.473
$ws.Range("I62").Value = 71561
$ws.Range("J62").Value = 11300
$ws.Range("K62").Value = 71561
$ws.Range("L62").Value = 11300
$ws.Range("M62").Value = -70937
$ws.Range("N62").Value = -12548
$ws.Range("H64").Value = 6200
$ws.Range("J64").Value = 4600
$ws.Range("L64").Value = 4600
$ws.Range("N64").Value = -5096
$ws.Range("H65").Value = 58874.473
$ws.Range("I65").Value = 71561
$ws.Range("J65").Value = 11300
$ws.Range("K65").Value = 357805
$ws.Range("L65").Value = 56500
$ws.Range("M65").Value = -354685
$ws.Range("N65").Value = -62740
$ws.Range("H67").Value = 6200
$ws.Range("J67").Value = 4600
$ws.Range("L67").Value = 4600
$ws.Range("N67").Value = -6316
$ws.Range("H70").Value = 1686.25
$ws.Range("I70").Value = 1050
$ws.Range("J70").Value = 2322.5
$ws.Range("K70").Value = 3150
$ws.Range("L70").Value = 6967.5
$ws.Range("M70").Value = -2880
$ws.Range("N70").Value = -7507.5
$ws.Range("H73").Value = 1686.25
$ws.Range("I73").Value = 1050
$ws.Range("J73").Value = 2322.5
$ws.Range("K73").Value = 3150
$ws.Range("L73").Value = 6967.5
$ws.Range("M73").Value = -2214
$ws.Range("N73").Value = -8839.5
$ws.Range("H112").Value = 3181.818
$ws.Range("J112").Value = 3635.2942
$ws.Range("L112").Value = 10905.8826
$ws.Range("N112").Value = -13121.8826
$ws.Range("H113").Value = 2983.9473
$ws.Range("I113").Value = 2886.6128
$ws.Range("J113").Value = 3415
$ws.Range("K113").Value = 2886.6128
$ws.Range("L113").Value = 3415
$ws.Range("M113").Value = 367.3872000000001
$ws.Range("N113").Value = -9923
$ws.Range("H138").Value = 3328.5
$ws.Range("I138").Value = 1691.4286
$ws.Range("J138").Value = 3763.671
$ws.Range("K138").Value = 5074.2858
$ws.Range("L138").Value = 11291.013
$ws.Range("M138").Value = 65.71420000000035
$ws.Range("N138").Value = -21571.013

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1318.129
$ws.Range("I74").Value = 1371.0869
$ws.Range("J74").Value = 1165.875
$ws.Range("K74").Value = 1371.0869
$ws.Range("L74").Value = 1165.875
$ws.Range("M74").Value = -497.0869
$ws.Range("N74").Value = -2913.875
$ws.Range("H77").Value = 1318.129
$ws.Range("I77").Value = 1371.0869
$ws.Range("J77").Value = 1165.875
$ws.Range("K77").Value = 6855.4345
$ws.Range("L77").Value = 5829.375
$ws.Range("M77").Value = -2487.4345
$ws.Range("N77").Value = -14565.375
$ws.Range("H132").Value = 1497.7833
$ws.Range("I132").Value = 1091.4255
$ws.Range("J132").Value = 2966.923
$ws.Range("K132").Value = 3274.2765
$ws.Range("L132").Value = 8900.769
$ws.Range("M132").Value = -744.2764999999999
$ws.Range("N132").Value = -13960.769
$ws.Range("H134").Value = 69245.8
$ws.Range("J134").Value = 69245.8
$ws.Range("L134").Value = 69245.8
$ws.Range("N134").Value = -79385.8
$ws.Range("H135").Value = 79714.5
$ws.Range("J135").Value = 79714.5
$ws.Range("L135").Value = 79714.5
$ws.Range("N135").Value = -89854.5
$ws.Range("H137").Value = 28254.5
$ws.Range("J137").Value = 47800
$ws.Range("L137").Value = 47800
$ws.Range("N137").Value = -58000
$ws.Range("H139").Value = 75800
$ws.Range("J139").Value = 75800
$ws.Range("L139").Value = 75800
$ws.Range("N139").Value = -86080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6306.8335
$ws.Range("I62").Value = 6554.7144
$ws.Range("J62").Value = 5959.8
$ws.Range("K62").Value = 6554.7144
$ws.Range("L62").Value = 5959.8
$ws.Range("M62").Value = -5930.7144
$ws.Range("N62").Value = -7207.8
$ws.Range("H65").Value = 6306.8335
$ws.Range("I65").Value = 6554.7144
$ws.Range("J65").Value = 5959.8
$ws.Range("K65").Value = 32773.572
$ws.Range("L65").Value = 29799
$ws.Range("M65").Value = -29653.572
$ws.Range("N65").Value = -36039
$ws.Range("H132").Value = 1537.4166
$ws.Range("I132").Value = 1094.3
$ws.Range("J132").Value = 3753
$ws.Range("K132").Value = 3282.9
$ws.Range("L132").Value = 11259
$ws.Range("M132").Value = -752.8999999999996
$ws.Range("N132").Value = -16319

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 70268.484
$ws.Range("J131").Value = 60721.53
$ws.Range("L131").Value = 182164.59
$ws.Range("N131").Value = -192244.59

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H70").Value = 4561.8
$ws.Range("I70").Value = 4266.6665
$ws.Range("J70").Value = 5004.5
$ws.Range("K70").Value = 4266.6665
$ws.Range("L70").Value = 5004.5
$ws.Range("M70").Value = -3996.6665
$ws.Range("N70").Value = -5544.5
$ws.Range("H73").Value = 4561.8
$ws.Range("I73").Value = 4266.6665
$ws.Range("J73").Value = 5004.5
$ws.Range("K73").Value = 4266.6665
$ws.Range("L73").Value = 5004.5
$ws.Range("M73").Value = -3330.6665
$ws.Range("N73").Value = -6876.5
$ws.Range("H122").Value = 1645717.9
$ws.Range("I122").Value = 2193965.5
$ws.Range("J122").Value = 975
$ws.Range("K122").Value = 6581896.5
$ws.Range("L122").Value = 2925
$ws.Range("M122").Value = -6579446.5
$ws.Range("N122").Value = -7825
$ws.Range("H132").Value = 3569.1282
$ws.Range("I132").Value = 3424.7856
$ws.Range("J132").Value = 3936.5454
$ws.Range("K132").Value = 10274.3568
$ws.Range("L132").Value = 11809.6362
$ws.Range("M132").Value = -7744.356800000001
$ws.Range("N132").Value = -16869.6362
$ws.Range("H133").Value = 66356
$ws.Range("J133").Value = 66356
$ws.Range("L133").Value = 66356
$ws.Range("N133").Value = -76476
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 722.9091
$ws.Range("I22").Value = 662.5
$ws.Range("J22").Value = 757.4286
$ws.Range("K22").Value = 662.5
$ws.Range("L22").Value = 757.4286
$ws.Range("M22").Value = -367.5
$ws.Range("N22").Value = -1347.4286
$ws.Range("H27").Value = 722.9091
$ws.Range("I27").Value = 662.5
$ws.Range("J27").Value = 757.4286
$ws.Range("K27").Value = 662.5
$ws.Range("L27").Value = 757.4286
$ws.Range("M27").Value = -555.5
$ws.Range("N27").Value = -971.4286
$ws.Range("H34").Value = 6999.5
$ws.Range("I34").Value = 6999
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 6999
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = -6827
$ws.Range("N34").Value = -7344
$ws.Range("H132").Value = 3193.93
$ws.Range("I132").Value = 2110.4138
$ws.Range("J132").Value = 4316.143
$ws.Range("K132").Value = 6331.241399999999
$ws.Range("L132").Value = 12948.429
$ws.Range("M132").Value = -3801.241399999999
$ws.Range("N132").Value = -18008.429
$ws.Range("H136").Value = 3857.7144
$ws.Range("I136").Value = 2151.973
$ws.Range("J136").Value = 9117.083000000001
$ws.Range("K136").Value = 6455.919
$ws.Range("L136").Value = 27351.249
$ws.Range("M136").Value = -3905.919
$ws.Range("N136").Value = -32451.249
